$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: the enquiry reply text is replaced with "0", and the enquiry
#     timestamp gets bumped (as if the entry was edited/re-saved). Briefly
#     formatting as Text keeps "0" stored as a string instead of letting
#     Excel auto-convert it to a number; resetting the style back to Normal
#     afterwards keeps the cell's formatting itself untouched. ---
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0"
$ws.Range("D4").Style = "Normal"

$ws.Range("F4").NumberFormat = "yyyy-MM-dd HH:mm:ss"
$ws.Range("F4").Value = 45768.54076372685

# --- Row 6: a reply "hi" is recorded, along with the reply timestamp ---
$ws.Range("E6").Value = "hi"
$ws.Range("F6").NumberFormat = "yyyy-MM-dd HH:mm:ss"
$ws.Range("F6").Value = 45767.863324502316
$ws.Range("G6").NumberFormat = "yyyy-MM-dd HH:mm:ss"
$ws.Range("G6").Value = 45768.54186474537

# --- Row 7: a brand-new enquiry is submitted ---
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "T2109876H"
$ws.Range("C7").Value = 1
$ws.Range("D7").Value = "hi"
$ws.Range("F7").NumberFormat = "yyyy-MM-dd HH:mm:ss"
$ws.Range("F7").Value = 45768.54314126157
